$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 1056:1058, pushing the existing data (and everything
# below it) down by three rows. This grows the used range from A1:T1146 to
# A1:T1149, matching the new weekly price report row added to the top of
# this block.
$ws.Rows("1056:1058").Insert()

# Populate the three newly-inserted rows with the new week's data (market /
# product metadata is identical to every other row in this report; only the
# date, quality grade, volume and price columns differ).
$newRows = @(1056, 1057, 1058)
$qualities = @("Especial", "Primera", "Segunda")

for ($i = 0; $i -lt 3; $i++) {
    $r = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 45013
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108002
    $ws.Cells.Item($r, 10).Value = "Mango"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $qualities[$i]
    $ws.Cells.Item($r, 13).Value = 576
    $ws.Cells.Item($r, 14).Value = 6000
    $ws.Cells.Item($r, 15).Value = 6500
    $ws.Cells.Item($r, 16).Value = 6250
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 19).Value = 1562
    $ws.Cells.Item($r, 20).Value = 4
}
